# Generate Report for Handoff
# - Refresh the "Latest Handoff"/"Latest HO Xliff Generate" timestamps for the
#   rows that just got a new handoff package, and stamp their Priority as "ht".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(8, 10, 11, 12, 13, 14)

foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-15 12:18:21"

    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-15 12:18:16"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-15 12:18:21"
}
